$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Row 3 / Row 4 -----------------------------------------------------
# New test data: user story id goes into B3, panel/template info into B4.
$ws.Range("B3").Value = "verifyBatteryStandbyAndAlarmLoadOnRbusAddDelete"
$ws.Range("B4").Value = "NGC-1928/T960 OR TC-71695 "

# Row 3 reverts to the default (non-custom) row height.
$ws.Rows.Item(3).AutoFit()
# Row 4 keeps its height but now flagged as an explicit custom height.
$ws.Rows.Item(4).RowHeight = 28.8

# --- Row 8 : panel / device data ---------------------------------------
$ws.Range("A8").Value = "FC64-2"

# C8 becomes an (quote-prefixed) empty string - use the leading apostrophe
# trick so it is stored as text (shared-string) rather than a truly blank
# cell, matching the "CPU Type" column convention used elsewhere (Sheet1
# column C uses the same empty-text-with-quote-prefix style).
$ws.Range("C8").Value = "'"

$ws.Range("F8").Value = 0.223
$ws.Range("G8").Value = 0.415

$ws.Range("N8").Value = "XIOM"
$ws.Range("P8").Value = "XIOM"

# Q8 / R8 use a quote-prefixed text style (s=10) in the source file even
# though they hold numbers; a plain .Value write downgrades that style
# because the new value no longer "needs" a quote prefix, so stash the
# original format in a scratch cell, write the value, then restore the
# format from the stash.
$ws.Range("Q8").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("Q8").Value = 0.049
$ws.Range("AB1").Copy()
$ws.Range("Q8").PasteSpecial(-4122)
$ws.Range("AB1").Clear()

$ws.Range("R8").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("R8").Value = 0.049
$ws.Range("AB1").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("AB1").Clear()

# S8 / T8 are the same story (quote-prefixed style s=19).
$ws.Range("S8").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("S8").Value = 0.294
$ws.Range("AB1").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("AB1").Clear()

$ws.Range("T8").Copy()
$ws.Range("AB1").PasteSpecial(-4122)
$ws.Range("T8").Value = 0.494
$ws.Range("AB1").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("AB1").Clear()

$ws.Range("U8").Value = 0.223
$ws.Range("V8").Value = 0.415

# --- Selection housekeeping ---------------------------------------------
# The saved file now has F8 selected (and Excel naturally stops forcing
# column I to the left edge once the selection moves back into view).
$ws.Activate()
$ws.Range("F8").Select()
